# Update the "想去人数" (want-to-go count) column F values as captured at a
# later scrape of the source site. Values only increase, reflecting new
# site visitors wanting to attend each event between commits.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibition)
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 93
$ws.Range("F8").Value = 2093
$ws.Range("F9").Value = 2006
$ws.Range("F10").Value = 1036
$ws.Range("F13").Value = 1626
$ws.Range("F17").Value = 76
$ws.Range("F18").Value = 96
$ws.Range("F20").Value = 545
$ws.Range("F21").Value = 652
$ws.Range("F22").Value = 336
$ws.Range("F23").Value = 11821
$ws.Range("F24").Value = 11840
$ws.Range("F27").Value = 1860
$ws.Range("F29").Value = 469

# Sheet "演出" (Performance)
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F6").Value = 36

# Sheet "本地生活" (Local life)
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 66

# Sheet "全部类型" (All types) - aggregates rows from the sheets above
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 66
$ws.Range("F6").Value = 93
$ws.Range("F10").Value = 2093
$ws.Range("F11").Value = 2006
$ws.Range("F12").Value = 1036
$ws.Range("F15").Value = 1626
$ws.Range("F20").Value = 76
$ws.Range("F22").Value = 96
$ws.Range("F24").Value = 545
$ws.Range("F25").Value = 652
$ws.Range("F26").Value = 336
$ws.Range("F27").Value = 11821
$ws.Range("F28").Value = 11840
$ws.Range("F31").Value = 1860
$ws.Range("F33").Value = 36
$ws.Range("F35").Value = 469
